$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reclassify ideologia (column B) using population-density based rule ---
# "Centro Dem" renamed to "Centro"; "Extinto" and "Centrao" folded into
# Direita/Centro/Esquerda per municipality
$ws.Range("B3").Value = "Centro"
$ws.Range("B6").Value = "Centro"
$ws.Range("B7").Value = "Centro"
$ws.Range("B11").Value = "Direita"
$ws.Range("B13").Value = "Centro"
$ws.Range("B14").Value = "Centro"
$ws.Range("B16").Value = "Centro"
$ws.Range("B17").Value = "Esquerda"
$ws.Range("B18").Value = "Centro"
$ws.Range("B19").Value = "Direita"
$ws.Range("B20").Value = "Centro"
$ws.Range("B21").Value = "Direita"
$ws.Range("B25").Value = "Centro"
$ws.Range("B26").Value = "Centro"
$ws.Range("B29").Value = "Esquerda"
$ws.Range("B31").Value = "Centro"
$ws.Range("B35").Value = "Centro"

# --- Update vote counts ---
$ws.Range("E2").Value = 8
$ws.Range("E3").Value = 131
$ws.Range("G3").Value = 125
$ws.Range("E5").Value = 503
$ws.Range("F5").Value = 286
$ws.Range("G5").Value = 272
$ws.Range("E6").Value = 1198
$ws.Range("F6").Value = 1021
$ws.Range("G6").Value = 1050
$ws.Range("E9").Value = 42
$ws.Range("G9").Value = 82
$ws.Range("E10").Value = 357
$ws.Range("F10").Value = 312
$ws.Range("G10").Value = 333
$ws.Range("E11").Value = 12
$ws.Range("G11").Value = 39
$ws.Range("E12").Value = 395
$ws.Range("F12").Value = 274
$ws.Range("G12").Value = 302
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = 12
$ws.Range("E16").Value = 544
$ws.Range("F16").Value = 473
$ws.Range("E22").Value = 12
$ws.Range("E23").Value = 306
$ws.Range("F23").Value = 435
$ws.Range("E24").Value = 57
$ws.Range("G24").Value = 86
$ws.Range("F25").Value = 493
$ws.Range("G25").Value = 532
$ws.Range("E26").Value = 797
$ws.Range("F26").Value = 706
$ws.Range("G26").Value = 807
$ws.Range("F27").Value = 22
$ws.Range("E30").Value = 555
$ws.Range("F30").Value = 642
$ws.Range("G30").Value = 255
$ws.Range("E31").Value = 402
$ws.Range("F31").Value = 290
$ws.Range("G31").Value = 259
$ws.Range("E33").Value = 74
$ws.Range("F33").Value = 100
$ws.Range("G33").Value = 100
